# Auto-generated edit script applying the diff to Alpha_Profits sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 92.947365
$ws.Range("I6").Value = 92.55556
$ws.Range("K6").Value = 277.66668
$ws.Range("M6").Value = -165.66668
$ws.Range("H8").Value = 307.1111
$ws.Range("J8").Value = 434.5
$ws.Range("L8").Value = 1303.5
$ws.Range("N8").Value = -1581.5
$ws.Range("H40").Value = 7039.8
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 7039.8
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 7039.8
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -7389.8
$ws.Range("H106").Value = 4185.3076
$ws.Range("I106").Value = 4426.5
$ws.Range("K106").Value = 4426.5
$ws.Range("M106").Value = -3795.5
$ws.Range("H116").Value = 60544.316
$ws.Range("J116").Value = 8726.272000000001
$ws.Range("L116").Value = 8726.272000000001
$ws.Range("N116").Value = -15610.272
$ws.Range("H125").Value = 2907.7273
$ws.Range("J125").Value = 2907.7273
$ws.Range("L125").Value = 26169.5457
$ws.Range("N125").Value = -31089.5457
$ws.Range("H135").Value = 907.8
$ws.Range("I135").Value = 652.6667
$ws.Range("J135").Value = 1928.3334
$ws.Range("K135").Value = 5874.0003
$ws.Range("L135").Value = 17355.0006
$ws.Range("M135").Value = -3339.0003
$ws.Range("N135").Value = -22425.0006

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3399.3438
$ws.Range("I32").Value = 3285.5862
$ws.Range("K32").Value = 3285.5862
$ws.Range("M32").Value = -2998.5862
$ws.Range("H97").Value = 522.5
$ws.Range("I97").Value = 295
$ws.Range("J97").Value = 750
$ws.Range("K97").Value = 295
$ws.Range("L97").Value = 750
$ws.Range("M97").Value = 201
$ws.Range("N97").Value = -1742
$ws.Range("H110").Value = 2896.75
$ws.Range("I110").Value = 2363.5908
$ws.Range("K110").Value = 2363.5908
$ws.Range("M110").Value = -318.5907999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2414.7778
$ws.Range("I99").Value = 2130.125
$ws.Range("K99").Value = 2130.125
$ws.Range("M99").Value = -632.125
$ws.Range("H105").Value = 5333.222
$ws.Range("I105").Value = 2998.1667
$ws.Range("K105").Value = 2998.1667
$ws.Range("M105").Value = -1251.1667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1496.9678
$ws.Range("J31").Value = 1197.9524
$ws.Range("L31").Value = 1197.9524
$ws.Range("N31").Value = -1787.9524
$ws.Range("H34").Value = 1496.9678
$ws.Range("J34").Value = 1197.9524
$ws.Range("L34").Value = 1197.9524
$ws.Range("N34").Value = -1601.9524
$ws.Range("H58").Value = 3114.7144
$ws.Range("I58").Value = 3016.0833
$ws.Range("K58").Value = 3016.0833
$ws.Range("M58").Value = -2813.0833
$ws.Range("H99").Value = 4369.143
$ws.Range("I99").Value = 3980.8333
$ws.Range("K99").Value = 3980.8333
$ws.Range("M99").Value = -2482.8333
$ws.Range("H126").Value = 4369.143
$ws.Range("I126").Value = 3980.8333
$ws.Range("K126").Value = 11942.4999
$ws.Range("M126").Value = -9472.499899999999
$ws.Range("H136").Value = 3114.7144
$ws.Range("I136").Value = 3016.0833
$ws.Range("K136").Value = 9048.249899999999
$ws.Range("M136").Value = -6498.249899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 360.33334
$ws.Range("I6").Value = 302.18182
$ws.Range("K6").Value = 906.54546
$ws.Range("M6").Value = -793.54546
$ws.Range("H8").Value = 514.625
$ws.Range("I8").Value = 514.625
$ws.Range("K8").Value = 1543.875
$ws.Range("M8").Value = -1404.875
$ws.Range("H12").Value = 519.1429000000001
$ws.Range("J12").Value = 518.6
$ws.Range("L12").Value = 1555.8
$ws.Range("N12").Value = -1901.8
$ws.Range("H63").Value = 377.5
$ws.Range("I63").Value = 246.875
$ws.Range("J63").Value = 900
$ws.Range("K63").Value = 740.625
$ws.Range("L63").Value = 2700
$ws.Range("M63").Value = 8.375
$ws.Range("N63").Value = -4198
$ws.Range("H64").Value = 6019.5835
$ws.Range("I64").Value = 3277.625
$ws.Range("J64").Value = 11503.5
$ws.Range("K64").Value = 9832.875
$ws.Range("L64").Value = 34510.5
$ws.Range("M64").Value = -9562.875
$ws.Range("N64").Value = -35050.5
$ws.Range("H66").Value = 377.5
$ws.Range("I66").Value = 246.875
$ws.Range("J66").Value = 900
$ws.Range("K66").Value = 2221.875
$ws.Range("L66").Value = 8100
$ws.Range("M66").Value = 1522.125
$ws.Range("N66").Value = -15588
$ws.Range("H67").Value = 6019.5835
$ws.Range("I67").Value = 3277.625
$ws.Range("J67").Value = 11503.5
$ws.Range("K67").Value = 9832.875
$ws.Range("L67").Value = 34510.5
$ws.Range("M67").Value = -8896.875
$ws.Range("N67").Value = -36382.5
$ws.Range("H70").Value = 2535.1667
$ws.Range("I70").Value = 1442
$ws.Range("K70").Value = 4326
$ws.Range("M70").Value = -4011
$ws.Range("H73").Value = 2535.1667
$ws.Range("I73").Value = 1442
$ws.Range("K73").Value = 4326
$ws.Range("M73").Value = -3234
$ws.Range("H113").Value = 761.1818
$ws.Range("J113").Value = 825.2857
$ws.Range("L113").Value = 2475.8571
$ws.Range("N113").Value = -6815.8571
$ws.Range("H131").Value = 53831.156
$ws.Range("I131").Value = 1043.6666
$ws.Range("J131").Value = 101339.9
$ws.Range("K131").Value = 3130.9998
$ws.Range("L131").Value = 304019.7
$ws.Range("M131").Value = 1909.0002
$ws.Range("N131").Value = -314099.7

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2746.4443
$ws.Range("I113").Value = 2784.5881
$ws.Range("J113").Value = 2098
$ws.Range("K113").Value = 2784.5881
$ws.Range("L113").Value = 2098
$ws.Range("M113").Value = -614.5880999999999
$ws.Range("N113").Value = -6438
$ws.Range("H122").Value = 3567.3684
$ws.Range("I122").Value = 4036.3076
$ws.Range("K122").Value = 12108.9228
$ws.Range("M122").Value = -9658.9228
$ws.Range("H126").Value = 3339.2104
$ws.Range("I126").Value = 2877.25
$ws.Range("J126").Value = 3462.4
$ws.Range("K126").Value = 8631.75
$ws.Range("L126").Value = 10387.2
$ws.Range("M126").Value = -6161.75
$ws.Range("N126").Value = -15327.2
$ws.Range("H132").Value = 4193
$ws.Range("I132").Value = 3515.4285
$ws.Range("J132").Value = 4870.5713
$ws.Range("K132").Value = 10546.2855
$ws.Range("L132").Value = 14611.7139
$ws.Range("M132").Value = -8016.2855
$ws.Range("N132").Value = -19671.7139

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4956.6
$ws.Range("I122").Value = 4720.75
$ws.Range("K122").Value = 14162.25
$ws.Range("M122").Value = -11712.25
$ws.Range("H132").Value = 2316
$ws.Range("I132").Value = 2040.6
$ws.Range("K132").Value = 6121.799999999999
$ws.Range("M132").Value = -3591.799999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3658
$ws.Range("I132").Value = 3322.5
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 9967.5
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -7437.5
